$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("adatok")

# Bump the "E" column (round-robin seat/position index) for rows 119-218:
# each block of rows shares a "D" value (the round number, 10..17) and the
# new E is the old E shifted by (D - 9) for that block.
for ($r = 119; $r -le 218; $r++) {
    $d = $ws.Cells.Item($r, 4).Value2
    $e = $ws.Cells.Item($r, 5).Value2
    $ws.Cells.Item($r, 5).Value = $e + ($d - 9)
}

# Make "adatok" the active sheet/tab with E5 selected (was "jatekos szinek"
# tab selected, with D7 selected there).
$ws.Activate()
$ws.Range("E5").Select()
